$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11: new entry - copy the date cell format from row 10 (A10) so we reuse
# the existing "date" cell style instead of creating a brand-new number format.
$ws.Range("A10").Copy($ws.Range("A11"))
$ws.Range("A11").Value = 44168
$ws.Range("B11").Value = 2
$ws.Range("D11").Value = "Passendes UI-Framework rausgesucht"

# Update the selected cell to match the author's final cursor position.
[void]$ws.Range("G18").Select()
